$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(72, 20, 12, 18, 30, 12)
    3  = @(12, 3, 2, 6, 1, 5)
    4  = @(20, 28, 20, 20, 20, 20)
    5  = @(2, 2, 2, 1, 4, 1)
    6  = @(4, 12, 5, 1, 10, 6)
    7  = @(5, 20, 8, 20, 13, 20)
    8  = @(20, 192, 342, 268, 25, 12)
    9  = @(15, 10, 12, 20, 17, 6)
    10 = @(36, 18, 21, 18, 22, 18)
    11 = @(32, 24, 16, 16, 32, 12)
    12 = @(8, 14, 14, 14, 10, 6)
    13 = @(36, 20, 18, 30, 42, 14)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
